$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 26.02.2022 12:00"

# Update D3 (Delta Cena) from text "+0.39" to numeric value 0.39
$ws.Range("D3").Value = 0.39

# Update E3 (Old Datum) from text date to numeric date serial value,
# carrying the same date/time number format as the other rows (E2, E4, ...).
$ws.Range("E3").Value = 44618.48967592593
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
